$wb = $excel.ActiveWorkbook

# --- Hoja1: update the "Conversion del dia" note with new rates ---
$ws1 = $wb.Worksheets.Item("Hoja1")
$ws1.Range("A1").Value = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 4.72 = 18702.83 pesos`n✅ 18702.83 pesos = 4.69 = 949.22 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

# --- tasas: update the tasas values in N10, O10, N12, O12 ---
$ws2 = $wb.Worksheets.Item("tasas")
$ws2.Range("N10").Value = 212
$ws2.Range("O10").Value = 3965
$ws2.Range("N12").Value = 3984
$ws2.Range("O12").Value = 202.2
